$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.727.30"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "3.122.00"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.16"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "625.26"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.16"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +8.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.378"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +4.10%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "3.124.65"
$ws.Range("E10").Value = "  -0.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.776"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +7.61%  "
$ws.Range("E12").Value = "  +3.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000256"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.57"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.32%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "91.679.63"
$ws.Range("E15").Value = "  +1.50%  "
$ws.Range("B16").Value = "Toncoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.50"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.49%  "
$ws.Range("D17").Value = "3.706.77"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").Value = "3.139.29"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.74"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("B20").Value = "PEPE"
$ws.Range("C20").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000222"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.31%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.78"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.87"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "447.93"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.12"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.88"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "92.09"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.96"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.92%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.250"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +25.77%  "
$ws.Range("E31").Value = "  +14.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.117"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +34.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.25"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.55%  "
$ws.Range("E34").Value = "  +24.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.168"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +10.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.69"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.53"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.15"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +23.09%  "
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.67"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.85%  "
$ws.Range("B40").Value = "PancakeSwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.92"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.43%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "492.06"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.30"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.422"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.16"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.93"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "156.46"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.696"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.10%  "
$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.58"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("B50").Value = "ImmutableX"
$ws.Range("C50").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.35"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.80"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.91%  "
